{"js": "// This edit fully restructures the document's content after the first\n// (heading) paragraph:\n//   - the one-sentence system overview becomes a labeled, multi-line block\n//     (Date / Time / OS Name / OS Version / Computer Name / IP Address),\n//   - the patch list gets explicit package versions/architectures,\n//   - the compliance, next-steps and risk-assessment sections are reworded\n//     and reformatted as \"*\"-bulleted lists,\n//   - several \"***\" section headers lose the space before their closing\n//     \"***\" (e.g. \"Patch Status Summary ***\" -> \"Patch Status Summary***\").\n//\n// Because nearly every paragraph after the heading is added, merged, or\n// reworded, the most reliable way to reproduce the target is to insert all\n// of the new paragraphs (as literal strings, in order) immediately before\n// the old content, then delete everything that used to be there.\n//\n// A trailing \"\\u000b\" (vertical tab / manual line break) in one of the\n// strings below represents a <w:br/> kept inside the same run as the text\n// before it, matching the target markup's <w:t>...</w:t><w:br/> pattern.\n// Each string becomes its own paragraph via insertParagraph.\n//\n// Inserting \"before\" the old, plain (\"Normal\"-styled, no <w:pPr>) \"System\n// Overview\" paragraph -- instead of \"after\" the Heading-2 title -- keeps\n// every newly inserted paragraph from inheriting the Heading 2 style.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// paragraphs.items[0] is the \"Operating System Patch Management RMF\n// Compliance\" heading; it is unchanged by this edit.\n// paragraphs.items[1] is the old \"*** System Overview ***\" paragraph --\n// it, and everything after it, is being replaced. It (and its successors)\n// are deleted below, after new content is inserted immediately before it.\nconst oldContentAnchor = paragraphs.items[1];\n\nconst newParagraphTexts = [\n  \"*** System Overview ***\\u000b\",\n  \"\\u000b\",\n  \"The system being monitored is a computer with the following specifications:\\u000b\",\n  \"Date: 04-06-2025\\u000b\",\n  \"Time: 16:02:09\\u000b\",\n  \"OS Name: kb322-18\\u000b\",\n  \"OS Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)\\u000b\",\n  \"Computer Name: kb322-18\\u000b\",\n  \"IP Address: 140.160.138.147\\u000b\",\n  \"\\u000b\",\n  \"*** Patch Status Summary***\\u000b\",\n  \"\\u000b\",\n  \"The following patches are pending and can be applied to bring the system up to date:\\u000b\",\n  \"* code/stable 1.99.0-1743632463 amd64\\u000b\",\n  \"* ure/stable-security 4:7.4.7-1+deb12u6 amd64\\u000b\",\n  \"* git-man/stable-security 1:2.39.5-0+deb12u2 all\\u000b\",\n  \"* git/stable-security 1:2.39.5-0+deb12u2 amd64\\u000b\",\n  \"\\u000b\",\n  \"These patches are related to security and can help protect the system against potential vulnerabilities.\\u000b\",\n  \"\\u000b\",\n  \"*** Compliance with RMF Controls ***\\u000b\",\n  \"\\u000b\",\n  \"To ensure compliance, it is recommended that:\\u000b\",\n  \"* The updates be reviewed and assessed for their impact on the system.\\u000b\",\n  \"* A corrective action plan be put in place to address any identified vulnerabilities or weaknesses.\\u000b\",\n  \"* Configuration management practices be followed to ensure that the system is properly secured and up to date.\\u000b\",\n  \"* Vulnerability checks be conducted regularly to identify potential security risks.\\u000b\",\n  \"\\u000b\",\n  \"*** Recommended next steps ***\\u000b\",\n  \"\\u000b\",\n  \"The recommended next steps are:\\u000b\",\n  \"* Review and assess the updates for their impact on the system.\\u000b\",\n  \"* Schedule patch deployments, if necessary.\\u000b\",\n  \"* Document the update process in a way that makes it easy to track and manage future patches.\\u000b\",\n  \"\\u000b\",\n  \"*** Risk Assessment***\\u000b\",\n  \"\\u000b\",\n  \"Based on the information provided, there is currently no indication of any pending security updates that could pose a risk to the system. However, regular vulnerability checks should be conducted to ensure that this remains the case.\"\n];\n\n// Insert every new paragraph directly before the (fixed) anchor -- each\n// insertParagraph(\"before\") lands just above the anchor, so looping in\n// forward order naturally reproduces the target order.\nfor (const text of newParagraphTexts) {\n  oldContentAnchor.insertParagraph(text, Word.InsertLocation.before);\n}\nawait context.sync();\n\n// Re-resolve the paragraph collection (indices shifted after the inserts\n// above) and delete the old anchor paragraph plus everything after it --\n// i.e. all of the stale, superseded original content.\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items\");\nawait context.sync();\n\nconst newCount = newParagraphTexts.length;\n// The old anchor is now at index newCount (1 heading + newCount inserted\n// paragraphs precede it); delete from the end back through that index.\nfor (let i = refreshed.items.length - 1; i >= newCount + 1; i--) {\n  refreshed.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# This edit fully restructures the document's content after the first\n# (heading) paragraph:\n#   - the one-sentence system overview becomes a labeled, multi-line block\n#     (Date / Time / OS Name / OS Version / Computer Name / IP Address),\n#   - the patch list gets explicit package versions/architectures,\n#   - the compliance, next-steps and risk-assessment sections are reworded\n#     and reformatted as \"*\"-bulleted lists,\n#   - several \"***\" section headers lose the space before their closing\n#     \"***\" (e.g. \"Patch Status Summary ***\" -> \"Patch Status Summary***\").\n#\n# Because nearly every paragraph after the heading is added, merged, or\n# reworded, the most reliable way to reproduce the target is to insert all\n# of the new paragraphs in one shot -- as a single block of text, in order,\n# using a vertical-tab ($VT / chr(11)) before each paragraph mark ($CR /\n# chr(13)) to represent a manual line break (<w:br/>) that stays inside the\n# same run as the preceding text, matching the target markup's\n# <w:t>...</w:t><w:br/> pattern -- immediately before the old content, and\n# then delete everything that used to be there.\n\n$d = $word.ActiveDocument\n\n$VT = [char]11\n$CR = [char]13\n\n# paragraph 1 is the \"Operating System Patch Management RMF Compliance\"\n# heading; it is unchanged by this edit.\n# paragraph 2 is the old \"*** System Overview ***\" paragraph -- it (and\n# every paragraph after it) is being replaced. It is used purely as an\n# anchor: the new content is inserted immediately before it, and then it\n# (plus everything after it) is deleted.\n$oldContentAnchor = $d.Paragraphs(2)\n\n$newContent = (\n  '*** System Overview ***' + $VT + $CR,\n  '' + $VT + $CR,\n  'The system being monitored is a computer with the following specifications:' + $VT + $CR,\n  'Date: 04-06-2025' + $VT + $CR,\n  'Time: 16:02:09' + $VT + $CR,\n  'OS Name: kb322-18' + $VT + $CR,\n  'OS Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)' + $VT + $CR,\n  'Computer Name: kb322-18' + $VT + $CR,\n  'IP Address: 140.160.138.147' + $VT + $CR,\n  '' + $VT + $CR,\n  '*** Patch Status Summary***' + $VT + $CR,\n  '' + $VT + $CR,\n  'The following patches are pending and can be applied to bring the system up to date:' + $VT + $CR,\n  '* code/stable 1.99.0-1743632463 amd64' + $VT + $CR,\n  '* ure/stable-security 4:7.4.7-1+deb12u6 amd64' + $VT + $CR,\n  '* git-man/stable-security 1:2.39.5-0+deb12u2 all' + $VT + $CR,\n  '* git/stable-security 1:2.39.5-0+deb12u2 amd64' + $VT + $CR,\n  '' + $VT + $CR,\n  'These patches are related to security and can help protect the system against potential vulnerabilities.' + $VT + $CR,\n  '' + $VT + $CR,\n  '*** Compliance with RMF Controls ***' + $VT + $CR,\n  '' + $VT + $CR,\n  'To ensure compliance, it is recommended that:' + $VT + $CR,\n  '* The updates be reviewed and assessed for their impact on the system.' + $VT + $CR,\n  '* A corrective action plan be put in place to address any identified vulnerabilities or weaknesses.' + $VT + $CR,\n  '* Configuration management practices be followed to ensure that the system is properly secured and up to date.' + $VT + $CR,\n  '* Vulnerability checks be conducted regularly to identify potential security risks.' + $VT + $CR,\n  '' + $VT + $CR,\n  '*** Recommended next steps ***' + $VT + $CR,\n  '' + $VT + $CR,\n  'The recommended next steps are:' + $VT + $CR,\n  '* Review and assess the updates for their impact on the system.' + $VT + $CR,\n  '* Schedule patch deployments, if necessary.' + $VT + $CR,\n  '* Document the update process in a way that makes it easy to track and manage future patches.' + $VT + $CR,\n  '' + $VT + $CR,\n  '*** Risk Assessment***' + $VT + $CR,\n  '' + $VT + $CR,\n  'Based on the information provided, there is currently no indication of any pending security updates that could pose a risk to the system. However, regular vulnerability checks should be conducted to ensure that this remains the case.' + $CR\n) -join ''\n\n$oldContentAnchor.Range.InsertBefore($newContent)\n\n# Delete the stale original content: after the insert above, the old\n# anchor paragraph (\"*** System Overview ***\") and everything that used to\n# follow it are still present, now pushed down after the newly inserted\n# paragraphs. Walk backwards from the last paragraph in the document down\n# to (and including) that old anchor, deleting each one.\n$newParaCount = 38\n$firstStaleIndex = 1 + $newParaCount + 1\nfor ($i = $d.Paragraphs.Count; $i -ge $firstStaleIndex; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
